$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Text-only edits (shared-string content changes, no style change)
# ------------------------------------------------------------------

# "...to my application..." -> "...to my ethical application..."
$ws.Range("A13").Value2 = "As a researcher, I want to be able to add any questionnaires I have designed to my ethical application so that they can get approval"

# "Shaun, Tim" -> "Shaun"
$ws.Range("F13").Value2 = "Shaun"

# "Laura + Michael" -> "Laura + Timothy" (independent of F15, which keeps "Laura + Michael")
$ws.Range("D20").Value2 = "Laura + Timothy"

# ------------------------------------------------------------------
# 2. F14: "Shanna" -> "Shanna, Tim (database needs setup first)" + wrap text on
# ------------------------------------------------------------------
$ws.Range("F14").Value2 = "Shanna, Tim (database needs setup first)"
$ws.Range("F14").WrapText = $true

# ------------------------------------------------------------------
# 3. Scrum table: snapshot formats BEFORE we overwrite the cells we copy
#    them from, so the captured look matches the pre-edit cell.
# ------------------------------------------------------------------

# A29 will become a new date row; it should look like A28 currently does
# (date-formatted) before A28's own format gets normalised to match A27.
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# A33 will receive the "M,B,S,I" caption that currently lives in A30;
# capture A30's current (plain) look first.
$ws.Range("A30").Copy()
$ws.Range("A33").PasteSpecial(-4122)

# Now normalise A28's date style to match A27's.
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)

# B28 should pick up the same (wrapped) look as the rest of row 27/28 text cells.
$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122)

# B29 (happiness score for the new day) should look like the plain style
# used elsewhere (e.g. F13) - no wrap.
$ws.Range("F13").Copy()
$ws.Range("B29").PasteSpecial(-4122)

# C29/E29/G29 (free text for the new day) should look like the wrapped
# text style used by B27/C27/E27/G27/C28/E28.
$ws.Range("C27").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("E27").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G27").Copy()
$ws.Range("G29").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4. Merge E29:F29 (new merge for the new day's "Issues" cell)
# ------------------------------------------------------------------
$ws.Range("E29:F29").Merge()

# ------------------------------------------------------------------
# 5. Fill in the new day-4 scrum entry (row 29) values
# ------------------------------------------------------------------
$ws.Range("A29").Value2 = 44224
$ws.Range("B29").Value2 = "5 b, 1 s, 1 i"
$ws.Range("C29").Value2 = "- Email code working!                           -MySQL stuff linked and working"
$ws.Range("E29").Value2 = "- overwhelmed (reassigned tasks to help)                                                  - some people can't work until other people have completed their work"
$ws.Range("G29").Value2 = "- get SQL reading                             -make questionairre pages working and branded                                          -add email stuff to send page"

# ------------------------------------------------------------------
# 6. Move the "M,B,S,I" legend caption from A30 down to A33, clearing A30
# ------------------------------------------------------------------
$ws.Range("A33").Value2 = "M,B,S,I"
$ws.Range("A30").ClearContents()
$ws.Range("A30").WrapText = $false
